# Rename the "NT WN_3,77" sheet to "NT WN_2_5" to reflect the new
# Teilvariante (3,77 -> 2,5) being simulated/evaluated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NT WN_3,77")
$ws.Name = "NT WN_2_5"

# Activate that sheet and update the scroll position / selection to match
# where the author left off editing (row ~19-24, cell J24 selected).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

# Update the driving input: Verdampferleistung (C22) changes from 3.77 kW
# to 2.5 kW. All the dependent formula cells (C8, D8, C24, C27, C29, D33,
# D34, D35, ...) recalc automatically off of this single input.
$ws.Range("C22").Value = 2.5

# Leave the selection where the author left it.
$ws.Range("J24").Select()
